# Update docs and re-generate reports
#
# The "Chart Report" sheet feeds six charts (whole-repo file-type counts,
# deployed/original schema version counts, deployed/original schema
# revision counts). The underlying counts were re-generated and a few of
# them changed; update the cached totals on the sheet so the linked
# charts pick up the new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart Report")

# Whole repository count of files grouped by type: xsd 270 -> 274
$ws.Range("B2").Value = 274

# Deployed Version count: 2019.1.000 135 -> 137
$ws.Range("B28").Value = 137

# Original Version count: 2019.1.000 135 -> 137
$ws.Range("B31").Value = 137

# Deployed Revision count: rev 20181201 128 -> 130
$ws.Range("B34").Value = 130

# Original Revision count: rev 20181201 128 -> 130
$ws.Range("B38").Value = 130

# Push the updated source data into the charts that reference these
# ranges so they stay in sync with the refreshed figures.
$chartCount = $ws.ChartObjects().Count
for ($i = 1; $i -le $chartCount; $i++) {
    $co = $ws.ChartObjects().Item($i)
    $chart = $co.Chart
    $chart.Refresh()
}
$wb.RefreshAll()
